$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 68
$ws.Range("B68").Value2 = 7609336.0
$ws.Range("F68").Value2 = "AlMuaidar"
$ws.Range("G68").Value2 = "AlWakrah SC"
$ws.Range("I68").Value2 = 4.0
$ws.Range("J68").Value2 = "A"
$ws.Range("K68").Value2 = 4.0
$ws.Range("L68").Value2 = 4.0
$ws.Range("M68").Value2 = 1.65
$ws.Range("N68").Value2 = 4.0
$ws.Range("P68").Value2 = 1.666
$ws.Range("Q68").Value2 = 0.75
$ws.Range("W68").Value2 = -1.0
$ws.Range("Y68").Value2 = 0.6659999999999999
$ws.Range("Z68").Value2 = -1.0
$ws.Range("AA68").Value2 = 0.875
$ws.Range("AB68").Value2 = 0.825
$ws.Range("AC68").Value2 = -1.0

# Row 69
$ws.Range("B69").Value2 = 7004626.0
$ws.Range("F69").Value2 = "Al Gharafa"
$ws.Range("G69").Value2 = "Qatar SC Doha"
$ws.Range("I69").Value2 = 1.0
$ws.Range("J69").Value2 = "H"
$ws.Range("K69").Value2 = 1.909
$ws.Range("L69").Value2 = 3.8
$ws.Range("M69").Value2 = 3.25
$ws.Range("N69").Value2 = 1.909
$ws.Range("P69").Value2 = 3.4
$ws.Range("Q69").Value2 = -0.5
$ws.Range("W69").Value2 = 0.909
$ws.Range("Y69").Value2 = -1.0
$ws.Range("Z69").Value2 = 0.925
$ws.Range("AA69").Value2 = -1.0
$ws.Range("AB69").Value2 = 0.0
$ws.Range("AC69").Value2 = -0.0

# Row 81
$ws.Range("B81").Value2 = 7840805.0
$ws.Range("F81").Value2 = "AlWakrah SC"
$ws.Range("G81").Value2 = "Al Markhiya"
$ws.Range("K81").Value2 = 1.062
$ws.Range("L81").Value2 = 11.0
$ws.Range("M81").Value2 = 17.0
$ws.Range("N81").Value2 = 1.363
$ws.Range("O81").Value2 = 4.75
$ws.Range("P81").Value2 = 7.0
$ws.Range("Q81").Value2 = -1.25
$ws.Range("R81").Value2 = 1.75
$ws.Range("S81").Value2 = 1.95
$ws.Range("U81").Value2 = 1.975
$ws.Range("V81").Value2 = 1.825
$ws.Range("Y81").Value2 = 6.0
$ws.Range("AA81").Value2 = 0.95
$ws.Range("AC81").Value2 = 0.4125

# Row 82
$ws.Range("B82").Value2 = 7840802.0
$ws.Range("F82").Value2 = "AlAhli Doha"
$ws.Range("G82").Value2 = "Umm Salal"
$ws.Range("K82").Value2 = 2.4
$ws.Range("L82").Value2 = 4.0
$ws.Range("M82").Value2 = 2.25
$ws.Range("N82").Value2 = 2.3
$ws.Range("O82").Value2 = 4.0
$ws.Range("P82").Value2 = 2.375
$ws.Range("Q82").Value2 = 0.0
$ws.Range("R82").Value2 = 1.875
$ws.Range("S82").Value2 = 1.925
$ws.Range("U82").Value2 = 2.0
$ws.Range("V82").Value2 = 1.8
$ws.Range("Y82").Value2 = 1.375
$ws.Range("AA82").Value2 = 0.925
$ws.Range("AC82").Value2 = 0.4

# Row 83
$ws.Range("B83").Value2 = 7840803.0
$ws.Range("F83").Value2 = "Al Sadd"
$ws.Range("G83").Value2 = "Qatar SC Doha"
$ws.Range("H83").Value2 = 3.0
$ws.Range("I83").Value2 = 0.0
$ws.Range("J83").Value2 = "H"
$ws.Range("K83").Value2 = 1.285
$ws.Range("L83").Value2 = 5.75
$ws.Range("M83").Value2 = 7.5
$ws.Range("N83").Value2 = 1.2
$ws.Range("O83").Value2 = 6.5
$ws.Range("P83").Value2 = 9.5
$ws.Range("Q83").Value2 = -2.0
$ws.Range("R83").Value2 = 1.95
$ws.Range("S83").Value2 = 1.85
$ws.Range("U83").Value2 = 1.95
$ws.Range("V83").Value2 = 1.85
$ws.Range("W83").Value2 = 0.2
$ws.Range("Y83").Value2 = -1.0
$ws.Range("Z83").Value2 = 0.95
$ws.Range("AA83").Value2 = -1.0
$ws.Range("AB83").Value2 = -1.0
$ws.Range("AC83").Value2 = 0.8500000000000001

# Row 84
$ws.Range("B84").Value2 = 7840685.0
$ws.Range("F84").Value2 = "Al Duhail"
$ws.Range("G84").Value2 = "Al Gharafa"
$ws.Range("H84").Value2 = 1.0
$ws.Range("I84").Value2 = 4.0
$ws.Range("J84").Value2 = "A"
$ws.Range("K84").Value2 = 2.1
$ws.Range("L84").Value2 = 3.8
$ws.Range("M84").Value2 = 2.9
$ws.Range("N84").Value2 = 1.8
$ws.Range("O84").Value2 = 4.0
$ws.Range("P84").Value2 = 3.6
$ws.Range("Q84").Value2 = -0.75
$ws.Range("R84").Value2 = 1.975
$ws.Range("S84").Value2 = 1.825
$ws.Range("U84").Value2 = 1.825
$ws.Range("V84").Value2 = 1.975
$ws.Range("W84").Value2 = -1.0
$ws.Range("Y84").Value2 = 2.6
$ws.Range("Z84").Value2 = -1.0
$ws.Range("AA84").Value2 = 0.825
$ws.Range("AB84").Value2 = 0.825
$ws.Range("AC84").Value2 = -1.0

# Row 90
$ws.Range("B90").Value2 = 7840809.0
$ws.Range("F90").Value2 = "Al Gharafa"
$ws.Range("G90").Value2 = "Al Sadd"
$ws.Range("H90").Value2 = 2.0
$ws.Range("I90").Value2 = 2.0
$ws.Range("J90").Value2 = "D"
$ws.Range("K90").Value2 = 5.0
$ws.Range("L90").Value2 = 4.75
$ws.Range("M90").Value2 = 1.45
$ws.Range("N90").Value2 = 5.25
$ws.Range("O90").Value2 = 5.0
$ws.Range("P90").Value2 = 1.4
$ws.Range("Q90").Value2 = 1.25
$ws.Range("R90").Value2 = 2.0
$ws.Range("S90").Value2 = 1.8
$ws.Range("T90").Value2 = 3.75
$ws.Range("U90").Value2 = 1.875
$ws.Range("V90").Value2 = 1.925
$ws.Range("W90").Value2 = -1.0
$ws.Range("X90").Value2 = 4.0
$ws.Range("Z90").Value2 = 1.0
$ws.Range("AB90").Value2 = 0.4375
$ws.Range("AC90").Value2 = -0.5

# Row 91
$ws.Range("B91").Value2 = 7840810.0
$ws.Range("F91").Value2 = "AlRayyan SC"
$ws.Range("G91").Value2 = "AlWakrah SC"
$ws.Range("H91").Value2 = 3.0
$ws.Range("I91").Value2 = 0.0
$ws.Range("J91").Value2 = "H"
$ws.Range("K91").Value2 = 2.0
$ws.Range("L91").Value2 = 3.6
$ws.Range("M91").Value2 = 3.1
$ws.Range("N91").Value2 = 2.15
$ws.Range("O91").Value2 = 3.4
$ws.Range("P91").Value2 = 2.9
$ws.Range("Q91").Value2 = -0.25
$ws.Range("R91").Value2 = 1.975
$ws.Range("S91").Value2 = 1.825
$ws.Range("T91").Value2 = 3.0
$ws.Range("U91").Value2 = 1.925
$ws.Range("V91").Value2 = 1.875
$ws.Range("W91").Value2 = 1.15
$ws.Range("X91").Value2 = -1.0
$ws.Range("Z91").Value2 = 0.9750000000000001
$ws.Range("AB91").Value2 = 0.0
$ws.Range("AC91").Value2 = -0.0

# Row 102
$ws.Range("B102").Value2 = 7840817.0
$ws.Range("F102").Value2 = "Al Markhiya"
$ws.Range("G102").Value2 = "Qatar SC Doha"
$ws.Range("H102").Value2 = 1.0
$ws.Range("I102").Value2 = 2.0
$ws.Range("J102").Value2 = "A"
$ws.Range("K102").Value2 = 4.0
$ws.Range("L102").Value2 = 3.6
$ws.Range("M102").Value2 = 1.75
$ws.Range("N102").Value2 = 3.5
$ws.Range("O102").Value2 = 3.75
$ws.Range("P102").Value2 = 1.833
$ws.Range("Q102").Value2 = 0.5
$ws.Range("R102").Value2 = 1.975
$ws.Range("S102").Value2 = 1.825
$ws.Range("T102").Value2 = 3.0
$ws.Range("U102").Value2 = 1.8
$ws.Range("V102").Value2 = 2.0
$ws.Range("W102").Value2 = -1.0
$ws.Range("Y102").Value2 = 0.833
$ws.Range("Z102").Value2 = -1.0
$ws.Range("AA102").Value2 = 0.825
$ws.Range("AB102").Value2 = 0.0
$ws.Range("AC102").Value2 = -0.0

# Row 103
$ws.Range("B103").Value2 = 7840818.0
$ws.Range("F103").Value2 = "Al Gharafa"
$ws.Range("G103").Value2 = "AlRayyan SC"
$ws.Range("H103").Value2 = 3.0
$ws.Range("I103").Value2 = 0.0
$ws.Range("J103").Value2 = "H"
$ws.Range("K103").Value2 = 2.5
$ws.Range("L103").Value2 = 3.4
$ws.Range("M103").Value2 = 2.5
$ws.Range("N103").Value2 = 2.3
$ws.Range("O103").Value2 = 3.5
$ws.Range("P103").Value2 = 2.7
$ws.Range("Q103").Value2 = 0.0
$ws.Range("R103").Value2 = 1.75
$ws.Range("S103").Value2 = 2.05
$ws.Range("T103").Value2 = 3.5
$ws.Range("U103").Value2 = 1.925
$ws.Range("V103").Value2 = 1.875
$ws.Range("W103").Value2 = 1.3
$ws.Range("Y103").Value2 = -1.0
$ws.Range("Z103").Value2 = 0.75
$ws.Range("AA103").Value2 = -1.0
$ws.Range("AB103").Value2 = -1.0
$ws.Range("AC103").Value2 = 0.875

# Row 108
$ws.Range("B108").Value2 = 7004655.0
$ws.Range("F108").Value2 = "AlRayyan SC"
$ws.Range("G108").Value2 = "Al Duhail"
$ws.Range("H108").Value2 = 2.0
$ws.Range("I108").Value2 = 0.0
$ws.Range("J108").Value2 = "H"
$ws.Range("K108").Value2 = 2.1
$ws.Range("L108").Value2 = 3.75
$ws.Range("M108").Value2 = 2.8
$ws.Range("N108").Value2 = 2.25
$ws.Range("O108").Value2 = 3.75
$ws.Range("P108").Value2 = 2.6
$ws.Range("Q108").Value2 = -0.25
$ws.Range("R108").Value2 = 1.975
$ws.Range("S108").Value2 = 1.825
$ws.Range("T108").Value2 = 3.5
$ws.Range("U108").Value2 = 1.925
$ws.Range("V108").Value2 = 1.775
$ws.Range("W108").Value2 = 1.25
$ws.Range("Y108").Value2 = -1.0
$ws.Range("Z108").Value2 = 0.9750000000000001
$ws.Range("AC108").Value2 = 0.7749999999999999

# Row 109
$ws.Range("B109").Value2 = 7882227.0
$ws.Range("F109").Value2 = "Al Markhiya"
$ws.Range("G109").Value2 = "Al Sadd"
$ws.Range("H109").Value2 = 1.0
$ws.Range("I109").Value2 = 2.0
$ws.Range("J109").Value2 = "A"
$ws.Range("K109").Value2 = 10.0
$ws.Range("L109").Value2 = 6.5
$ws.Range("M109").Value2 = 1.2
$ws.Range("N109").Value2 = 11.0
$ws.Range("O109").Value2 = 7.5
$ws.Range("P109").Value2 = 1.142
$ws.Range("Q109").Value2 = 2.25
$ws.Range("R109").Value2 = 1.95
$ws.Range("S109").Value2 = 1.85
$ws.Range("T109").Value2 = 3.75
$ws.Range("U109").Value2 = 1.825
$ws.Range("V109").Value2 = 1.975
$ws.Range("W109").Value2 = -1.0
$ws.Range("Y109").Value2 = 0.1419999999999999
$ws.Range("Z109").Value2 = 0.95
$ws.Range("AC109").Value2 = 0.9750000000000001

# Row 110
$ws.Range("N110").Value2 = 5.0
$ws.Range("P110").Value2 = 1.5
$ws.Range("R110").Value2 = 2.0
$ws.Range("S110").Value2 = 1.8
$ws.Range("U110").Value2 = 1.9
$ws.Range("V110").Value2 = 1.9

# Row 111
$ws.Range("N111").Value2 = 2.15
$ws.Range("P111").Value2 = 2.9
$ws.Range("R111").Value2 = 1.9
$ws.Range("S111").Value2 = 1.9
$ws.Range("T111").Value2 = 3.25
$ws.Range("U111").Value2 = 2.025
$ws.Range("V111").Value2 = 1.775
